$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 66
$ws1.Range("F4").Value = 1487
$ws1.Range("F5").Value = 576
$ws1.Range("F7").Value = 11007
$ws1.Range("F8").Value = 11007
$ws1.Range("F11").Value = 314
$ws1.Range("F13").Value = 751
$ws1.Range("F14").Value = 12204
$ws1.Range("F15").Value = 12724
$ws1.Range("F22").Value = 30

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 66
$ws4.Range("F5").Value = 1487
$ws4.Range("F6").Value = 576
$ws4.Range("F8").Value = 11007
$ws4.Range("F9").Value = 11007
$ws4.Range("F12").Value = 314
$ws4.Range("F14").Value = 751
$ws4.Range("F15").Value = 12204
$ws4.Range("F16").Value = 12724
$ws4.Range("F23").Value = 30
